$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B18").Value = 0.0014449000009335499
$ws.Range("C18").Value = 0.386416800000006

$ws.Range("B18").Select()
